# Auto-generated script to apply row data re-ordering per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 16895948
$ws.Range("B2").Value = 89406
$ws.Range("D2").Value = 'NT'
$ws.Range("E2").Value = 1204
$ws.Range("F2").Value = 'Gränsticka'
$ws.Range("G2").Value = 'Phellopilus nigrolimitatus'
$ws.Range("H2").Value = '(Romell) Niemelä, T.Wagner & M.Fisch.'
$ws.Range("P2").Value = 'Kycklingvattnet, Ö om, 300 m NO Sätertjärnen, Jmt'
$ws.Range("Q2").Value = 461487.2423814723
$ws.Range("R2").Value = 7164776.231211329
$ws.Range("AI2").Value = 'Storluckig gammal, f.d. betad fjällgranskog'
$ws.Range("AN2").Value = 1
$ws.Range("AO2").Value = '1 substratenheter # grov hård granlåga'

# Row 3
$ws.Range("A3").Value = 16895952
$ws.Range("B3").Value = 77668
$ws.Range("D3").Value = 'VU'
$ws.Range("E3").Value = 1249
$ws.Range("F3").Value = 'Norsk näverlav'
$ws.Range("G3").Value = 'Platismatia norvegica'
$ws.Range("H3").Value = '(Lynge) W.L.Culb. & C.F.Culb.'
$ws.Range("P3").Value = 'Kycklingvattnet, Ö om, 400 m NO Sätertjärnen, Jmt'
$ws.Range("Q3").Value = 461636.9326544968
$ws.Range("R3").Value = 7164647.242933135
$ws.Range("AI3").Value = 'Storluckig gammal, f.d. betad fjällgranskog'
$ws.Range("AN3").Value = 5
$ws.Range("AO3").Value = '5 substratenheter # grenar av senvuxen gammal klonbildande gran'

# Row 4
$ws.Range("A4").Value = 16895947
$ws.Range("B4").Value = 77668
$ws.Range("D4").Value = 'VU'
$ws.Range("E4").Value = 1249
$ws.Range("F4").Value = 'Norsk näverlav'
$ws.Range("G4").Value = 'Platismatia norvegica'
$ws.Range("H4").Value = '(Lynge) W.L.Culb. & C.F.Culb.'
$ws.Range("P4").Value = 'Kycklingvattnet, Ö om, 200 m NO Sätertjärnen, Jmt'
$ws.Range("Q4").Value = 461433.4034202741
$ws.Range("R4").Value = 7164700.488456889
$ws.Range("AI4").Value = 'Storluckig gammal, f.d. betad fjällgranskog'
$ws.Range("AN4").Value = 3
$ws.Range("AO4").Value = '3 substratenheter # grenar på mkt grova granar och torrgran'

# Row 5
$ws.Range("A5").Value = 16895956
$ws.Range("B5").Value = 89406
$ws.Range("D5").Value = 'NT'
$ws.Range("E5").Value = 1204
$ws.Range("F5").Value = 'Gränsticka'
$ws.Range("G5").Value = 'Phellopilus nigrolimitatus'
$ws.Range("H5").Value = '(Romell) Niemelä, T.Wagner & M.Fisch.'
$ws.Range("P5").Value = 'Kycklingvattnet, Ö om, ca 400 m NO Sätertjärnen, Jmt'
$ws.Range("Q5").Value = 461898.1476932919
$ws.Range("R5").Value = 7164702.346384764
$ws.Range("AI5").Value = 'Storluckig gammal, f.d. betad fjällbjörkskog'
$ws.Range("AN5").Value = 1
$ws.Range("AO5").Value = '1 substratenheter # grov ngt murken granlåga'

# Row 6
$ws.Range("A6").Value = 16895951
$ws.Range("B6").Value = 77668
$ws.Range("D6").Value = 'VU'
$ws.Range("E6").Value = 1249
$ws.Range("F6").Value = 'Norsk näverlav'
$ws.Range("G6").Value = 'Platismatia norvegica'
$ws.Range("H6").Value = '(Lynge) W.L.Culb. & C.F.Culb.'
$ws.Range("P6").Value = 'Kycklingvattnet, Ö om, 400 m NO Sätertjärnen, Jmt'
$ws.Range("Q6").Value = 461560.5512648276
$ws.Range("R6").Value = 7164760.27418993
$ws.Range("AI6").Value = 'Storluckig gammal, f.d. betad fjällbjörkskog'
$ws.Range("AN6").Value = 1
$ws.Range("AO6").Value = '1 substratenheter # stammen på äldre sälg'

# Row 7
$ws.Range("A7").Value = 16895946
$ws.Range("B7").Value = 77668
$ws.Range("D7").Value = 'VU'
$ws.Range("E7").Value = 1249
$ws.Range("F7").Value = 'Norsk näverlav'
$ws.Range("G7").Value = 'Platismatia norvegica'
$ws.Range("H7").Value = '(Lynge) W.L.Culb. & C.F.Culb.'
$ws.Range("P7").Value = 'Kycklingvattnet, Ö om, 400 m NV Sätertjärnen, Jmt'
$ws.Range("Q7").Value = 461019.1146955636
$ws.Range("R7").Value = 7164878.383126226
$ws.Range("AI7").Value = 'Storluckig gammal, f.d. betad fjällgranskog'
$ws.Range("AN7").Value = 1
$ws.Range("AO7").Value = '1 substratenheter # rikligt på grenar och stam av mkt gammal gran'

# Row 8
$ws.Range("A8").Value = 16895955
$ws.Range("B8").Value = 89406
$ws.Range("D8").Value = 'NT'
$ws.Range("E8").Value = 1204
$ws.Range("F8").Value = 'Gränsticka'
$ws.Range("G8").Value = 'Phellopilus nigrolimitatus'
$ws.Range("H8").Value = '(Romell) Niemelä, T.Wagner & M.Fisch.'
$ws.Range("P8").Value = 'Kycklingvattnet, Ö om, 500 m NO Sätertjärnen, Jmt'
$ws.Range("Q8").Value = 461817.4838700104
$ws.Range("R8").Value = 7164851.487865292
$ws.Range("AI8").Value = 'Storluckig gammal, f.d. betad fjällbjörkskog'
$ws.Range("AN8").Value = 1
$ws.Range("AO8").Value = '1 substratenheter # grov ngt murken granlåga'

# Row 9
$ws.Range("A9").Value = 16895945
$ws.Range("B9").Value = 77668
$ws.Range("D9").Value = 'VU'
$ws.Range("E9").Value = 1249
$ws.Range("F9").Value = 'Norsk näverlav'
$ws.Range("G9").Value = 'Platismatia norvegica'
$ws.Range("H9").Value = '(Lynge) W.L.Culb. & C.F.Culb.'
$ws.Range("P9").Value = 'Kycklingvattnet, Ö om, 500 m NV Sätertjärnen, Jmt'
$ws.Range("Q9").Value = 461022.9782998873
$ws.Range("R9").Value = 7164978.369098279
$ws.Range("AI9").Value = 'Storluckig gammal, f.d. betad fjällgranskog'
$ws.Range("AN9").Value = 1
$ws.Range("AO9").Value = '1 substratenheter # rikligt på grova grangrenar'

# Row 10
$ws.Range("A10").Value = 16895950
$ws.Range("B10").Value = 76862
$ws.Range("D10").Value = 'LC'
$ws.Range("E10").Value = 6443
$ws.Range("F10").Value = 'Sotlav'
$ws.Range("G10").Value = 'Acolium inquinans'
$ws.Range("H10").Value = '(Sm.) A.Massal.'
$ws.Range("P10").Value = 'Kycklingvattnet, Ö om, 300 m NO Sätertjärnen, Jmt'
$ws.Range("Q10").Value = 461477.6555552008
$ws.Range("R10").Value = 7164766.04831084
$ws.Range("AI10").Value = 'Storluckig gammal, f.d. betad fjällgranskog'
$ws.Range("AN10").Value = 1
$ws.Range("AO10").Value = '1 substratenheter # torrgrenar på gammal gran i granklon'

# Row 11
$ws.Range("A11").Value = 16895959
$ws.Range("B11").Value = 77668
$ws.Range("D11").Value = 'VU'
$ws.Range("E11").Value = 1249
$ws.Range("F11").Value = 'Norsk näverlav'
$ws.Range("G11").Value = 'Platismatia norvegica'
$ws.Range("H11").Value = '(Lynge) W.L.Culb. & C.F.Culb.'
$ws.Range("P11").Value = 'Kycklingvattnet, Ö om, ca 1 km Ö Sätertjärnen, Jmt'
$ws.Range("Q11").Value = 462630.2221353759
$ws.Range("R11").Value = 7164280.166798776
$ws.Range("AI11").Value = 'Storluckig, gammal fjällgranskog på myrholme'
$ws.Range("AN11").Value = 10
$ws.Range("AO11").Value = '10 substratenheter # grenar och stammar av gammal klonbildande gran'

# Row 12
$ws.Range("A12").Value = 16895958
$ws.Range("B12").Value = 77668
$ws.Range("D12").Value = 'VU'
$ws.Range("E12").Value = 1249
$ws.Range("F12").Value = 'Norsk näverlav'
$ws.Range("G12").Value = 'Platismatia norvegica'
$ws.Range("H12").Value = '(Lynge) W.L.Culb. & C.F.Culb.'
$ws.Range("P12").Value = 'Kycklingvattnet, Ö om, ca 1 km Ö Sätertjärnen, Jmt'
$ws.Range("Q12").Value = 462429.0675681746
$ws.Range("R12").Value = 7164389.574527718
$ws.Range("AI12").Value = 'Storluckig, gammal fjällgranskog på myrholme'
$ws.Range("AN12").Value = 4
$ws.Range("AO12").Value = '4 substratenheter # grenar och stammar av gammal klonbildande gran'

# Row 13
$ws.Range("A13").Value = 16895965
$ws.Range("B13").Value = 77668
$ws.Range("D13").Value = 'VU'
$ws.Range("E13").Value = 1249
$ws.Range("F13").Value = 'Norsk näverlav'
$ws.Range("G13").Value = 'Platismatia norvegica'
$ws.Range("H13").Value = '(Lynge) W.L.Culb. & C.F.Culb.'
$ws.Range("P13").Value = 'Kycklingvattnet, Ö om, ca 600 m SSO Sätertjärnen, Jmt'
$ws.Range("Q13").Value = 461542.5122580806
$ws.Range("R13").Value = 7163913.259641338
$ws.Range("AI13").Value = 'Skiktad, gammal fjällgranskog'
$ws.Range("AN13").Value = 2
$ws.Range("AO13").Value = '2 substratenheter # grenar av gammal torrgran'

